# "Ui Test 3rd Sprint" report - update the heading of the 3rd-sprint section
# and fold the now-obsolete "τρίτης φάσης Ui:" wording into a single,
# shorter run ("διεπαφής (Ανανεωμένα μενού)").

$d = $word.ActiveDocument

# 1) "Έναρξη " -> "Τρίτο μέρος "  (first run of the title paragraph)
$d.Content.Find.Execute("Έναρξη ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Τρίτο μέρος ", 2) | Out-Null

# 2) "test" -> "test "  (second run, now needs a trailing space since the
#    following runs that used to supply it are being collapsed away)
$d.Content.Find.Execute("test", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "test ", 2) | Out-Null

# 3) Collapse " τρίτης φάσης Ui:" (originally five separate runs, two of
#    them wrapped in spell-check proofErr markers) into one new run reading
#    "διεπαφής (Ανανεωμένα μενού)".
$d.Content.Find.Execute(" τρίτης φάσης Ui:", $true, $false, $false, $false, `
                         $false, $true, 1, $false, `
                         "διεπαφής (Ανανεωμένα μενού)", 2) | Out-Null
